# Build script: fixes the LOQ4052 syllabus sheet.
# 1) Inserts a new row 13 (shifts everything below down by one) to hold the
#    "5840560 - Marco Antonio Carvalho Pereira" (professor) entry that was
#    previously mis-placed one row too high (it lived under "Objetivos:" and
#    again under "Metodo:").
# 2) Fixes the resulting/pre-existing off-by-one content so every label in
#    column A lines up with its correct description in columns B/C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 13 - pushes old rows 13-23 down to 14-24
# and carries each row's height/formatting along with it.
$ws.Rows("13").Insert()

# New row 13: professor responsible for the course.
$ws.Range("B13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C13").Value = "5840560 - Marco Antonio Carvalho Pereira"

# The freshly-inserted row only carries column A's formatting by default;
# pull the real column B/C formatting (wrap text, red font, ...) down from
# the row underneath, then drop the leftover empty A13 cell - row 13 has no
# label in column A.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("A13").Clear()

# Row 10 ("Objetivos:") previously (incorrectly) held the professor's name;
# it should hold the actual course objectives text.
$ws.Range("B10").Value = "Complementar a formação multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, tópicos atuais e relevantes sobre gestão e produção."
$ws.Range("C10").Value = "Complementar a formação multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, tópicos atuais e relevantes sobre gestão e produção."

# Row 14 ("Programa resumido:") should describe the topic, not "Semestral".
$ws.Range("B14").Value = "A definir, de acordo com o tópico programado."
$ws.Range("C14").Value = "A definir, de acordo com o tópico programado."

# Row 16 ("Programa:") should hold the actual syllabus content, not a date.
$ws.Range("B16").Value = "O conteúdo desta disciplina será de acordo com o tópico a ser programado, devendo abordar assuntos complementares a formação de um profissional de Engenharia."
$ws.Range("C16").Value = "O conteúdo desta disciplina será de acordo com o tópico a ser programado, devendo abordar assuntos complementares a formação de um profissional de Engenharia."

# Row 19 ("Método:") should hold the teaching method description.
$ws.Range("B19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("C19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."

# Row 20 ("Critério:") should hold the grading criteria.
$ws.Range("B20").Value = "Provas e trabalhos"
$ws.Range("C20").Value = "Provas e trabalhos"

# Row 21 ("Norma de recuperação:") should hold the recovery-exam rule.
$ws.Range("B21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."

# Row 22 ("Bibliografia:") should hold the bibliography text.
$ws.Range("B22").Value = "Textos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas na área de gestão e produção."
$ws.Range("C22").Value = "Textos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas na área de gestão e produção."
